$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the "Emilia" row (row 3: 005535788 | Emilia | 21000.56)
$ws.Rows.Item(3).Delete()

# 2) Insert a new row before the "Nabor" row (originally row 30, now row 29
#    after the deletion above) and fill it with Jo's new data.
$ws.Rows.Item(29).Insert()
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "005324981"
$ws.Cells.Item(29, 2).Value = "Jo"
$ws.Cells.Item(29, 3).Value = 495.94

# 3) Delete the old "Jo" row (005324981 | Jo | -16467.4) and the
#    "Venia" row (004813166 | Venia | -16979.46) that followed it.
#    Net row-count change so far is zero (one delete + one insert),
#    so these are still at their original absolute row numbers: 468, 469.
$ws.Rows.Item(468).Delete()
$ws.Rows.Item(468).Delete()
